$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same rows of data and both
# need their "想去人数" (F column) figures refreshed.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 57
    $ws.Range("F3").Value = 342
    $ws.Range("F5").Value = 94
}
